# Update "想去人数" (column F) figures on both the "展览" and "全部类型"
# worksheets to reflect the regenerated scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> list of (cell, newValue) updates.
$updates = @{
    "展览"     = @(
        @{ Cell = "F3";  Value = 559 },
        @{ Cell = "F9";  Value = 1158 },
        @{ Cell = "F10"; Value = 16283 },
        @{ Cell = "F14"; Value = 6351 },
        @{ Cell = "F15"; Value = 636 },
        @{ Cell = "F28"; Value = 890 },
        @{ Cell = "F30"; Value = 5047 },
        @{ Cell = "F39"; Value = 73 }
    )
    "全部类型" = @(
        @{ Cell = "F3";  Value = 559 },
        @{ Cell = "F9";  Value = 1158 },
        @{ Cell = "F10"; Value = 16283 },
        @{ Cell = "F14"; Value = 6351 },
        @{ Cell = "F15"; Value = 636 },
        @{ Cell = "F28"; Value = 890 },
        @{ Cell = "F30"; Value = 5047 },
        @{ Cell = "F40"; Value = 73 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates[$sheetName]) {
        $ws.Range($u.Cell).Value = $u.Value
    }
}

$wb.Save()
